$wb = $excel.ActiveWorkbook

# Sheet "展览" - update 想去人数 (interested count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 291
$ws1.Range("F4").Value = 7836
$ws1.Range("F5").Value = 5729
$ws1.Range("F6").Value = 475
$ws1.Range("F11").Value = 286
$ws1.Range("F12").Value = 60

# Sheet "全部类型" - mirrors same rows, with the last two entries shifted to F13/F14
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 291
$ws4.Range("F4").Value = 7836
$ws4.Range("F5").Value = 5729
$ws4.Range("F6").Value = 475
$ws4.Range("F13").Value = 286
$ws4.Range("F14").Value = 60
